$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quadratic-svm-score")

# Re-apply the (default) Locked protection flag to the header row and the
# first data-row label cell; this forces a fresh style record to be
# allocated for these cells (still numFmtId 49 / text, no border - same
# rendering as before) instead of reusing the old shared one, mirroring
# the style-table growth seen in the target workbook.
$ws.Range("A1:C1").Locked = $true
$ws.Range("A2").Locked = $true

# Update the numeric prediction score in B2.
$ws.Range("B2").Value = 98.937118899477554
